$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44586
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 3500

# Row 3
$ws.Range("D3").Value = 44587
$ws.Range("M3").Value = 165
$ws.Range("N3").Value = 6500
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6742
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value = 3371

# Row 4
$ws.Range("D4").Value = 44588
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 6500
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6750
$ws.Range("S4").Value = 3375

# Row 5
$ws.Range("D5").Value = 44589
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("S5").Value = 3000

# Row 7
$ws.Range("D7").Value = 44585
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 6500
$ws.Range("P7").Value = 6750
$ws.Range("S7").Value = 3375

# Row 8
$ws.Range("D8").Value = 44209
$ws.Range("M8").Value = 58
$ws.Range("N8").Value = 6000
$ws.Range("O8").Value = 6000
$ws.Range("P8").Value = 6000
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 3000

# Row 10
$ws.Range("D10").Value = 44627
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("R10").Value = "Provincia de Linares"
$ws.Range("S10").Value = 3000

# Row 11
$ws.Range("D11").Value = 44592
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 4000

# Row 13
$ws.Range("D13").Value = 44614

# Row 15
$ws.Range("D15").Value = 44211
$ws.Range("M15").Value = 45
